$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header last-install-check date (O1) ---
$ws.Range("O1").Value = 45519

# --- Remove all existing hyperlinks (will be rebuilt after data is rewritten) ---
$ws.Hyperlinks.Delete()

# --- Copy date-cell (C/D) number formatting down into the newly added rows 30-34 ---
# (columns C and D have no sheet-level default style, so fresh cells need the date format copied explicitly)
$ws.Range("C2").Copy($ws.Range("C30")) | Out-Null
$ws.Range("D2").Copy($ws.Range("D30")) | Out-Null
$ws.Range("C2").Copy($ws.Range("C31")) | Out-Null
$ws.Range("D2").Copy($ws.Range("D31")) | Out-Null
$ws.Range("C2").Copy($ws.Range("C32")) | Out-Null
$ws.Range("D2").Copy($ws.Range("D32")) | Out-Null
$ws.Range("C2").Copy($ws.Range("C33")) | Out-Null
$ws.Range("D2").Copy($ws.Range("D33")) | Out-Null
$ws.Range("C2").Copy($ws.Range("C34")) | Out-Null
$ws.Range("D2").Copy($ws.Range("D34")) | Out-Null

# --- Write World / Version / Release Date / Install Date / Source for each row ---
# Row 9 (Final Fantasy 5 Career Day) is untouched by this update, so it is skipped.

# Row 2: Animal Well
$ws.Range("A2").Value = "Animal Well"
$ws.Range("B2").Value = "0.4.0"
$ws.Range("C2").Value = 45529
$ws.Range("D2").Value = 45533
$ws.Range("E2").Value = "https://github.com/ScipioWright/Archipelago-SW/releases"

# Row 3: A Robot Named Fight!
$ws.Range("A3").Value = "A Robot Named Fight!"
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = 45274
$ws.Range("D3").Value = 45422
$ws.Range("E3").Value = "https://discord.com/channels/731205301247803413/1169389087371841708"

# Row 4: Celeste
$ws.Range("A4").Value = "Celeste"
$ws.Range("B4").Value = "0.3.0"
$ws.Range("C4").Value = 45308
$ws.Range("D4").Value = 45410
$ws.Range("E4").Value = "https://github.com/doshyw/CelesteArchipelago/releases"

# Row 5: Chrono Trigger Jets of Time
$ws.Range("A5").Value = "Chrono Trigger Jets of Time"
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = 45519
$ws.Range("E5").Value = "https://wiki.ctjot.com/doku.php?id=multiworld"

# Row 6: CrossCode
$ws.Range("A6").Value = "CrossCode"
$ws.Range("B6").Value = "0.6.1"
$ws.Range("C6").Value = 45533
$ws.Range("D6").Value = 45533
$ws.Range("E6").Value = "https://github.com/CodeTriangle/CCMultiworldRandomizer/releases"

# Row 7: EarthBound
$ws.Range("A7").Value = "EarthBound"
$ws.Range("B7").Value = "2.0.1"
$ws.Range("C7").Value = 45529
$ws.Range("D7").Value = 45533
$ws.Range("E7").Value = "https://github.com/PinkSwitch/Archipelago/releases/"

# Row 8: Final Fantasy 12 Open World
$ws.Range("A8").Value = "Final Fantasy 12 Open World"
$ws.Range("B8").Value = "0.3.8"
$ws.Range("C8").Value = 45528
$ws.Range("D8").Value = 45533
$ws.Range("E8").Value = "https://github.com/Bartz24/Archipelago/releases"

# Row 10: Final Fantasy 6 Worlds Collide
$ws.Range("A10").Value = "Final Fantasy 6 Worlds Collide"
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = 45467
$ws.Range("D10").Value = 45478
$ws.Range("E10").Value = "https://discord.com/channels/731205301247803413/1022545979146252288"

# Row 11: Fire Emblem 8
$ws.Range("A11").Value = "Fire Emblem 8"
$ws.Range("B11").Value = "alpha-0.3.2"
$ws.Range("C11").Value = 45425
$ws.Range("D11").Value = 45478
$ws.Range("E11").Value = "https://github.com/CT075/Archipelago/releases"

# Row 12: Inscryption
$ws.Range("A12").Value = "Inscryption"
$ws.Range("B12").Value = "0.3.0"
$ws.Range("C12").Value = 45479
$ws.Range("D12").Value = 45507
$ws.Range("E12").Value = "https://github.com/DrBibop/Archipelago_Inscryption/releases"

# Row 13: Jak and Daxter
$ws.Range("A13").Value = "Jak and Daxter"
$ws.Range("B13").Value = "0.0.10"
$ws.Range("C13").Value = 45528
$ws.Range("D13").Value = 45533
$ws.Range("E13").Value = "https://github.com/ArchipelaGOAL/Archipelago/releases"

# Row 14: Kingdom Hearts: Birth by Sleep
$ws.Range("A14").Value = "Kingdom Hearts: Birth by Sleep"
$ws.Range("B14").Value = "0.0.3"
$ws.Range("C14").Value = 45504
$ws.Range("D14").Value = 45507
$ws.Range("E14").Value = "https://github.com/gaithernOrg/ArchipelagoKHBBS/releases/"

# Row 15: Kingdom Hearts 1
$ws.Range("A15").Value = "Kingdom Hearts 1"
$ws.Range("B15").Value = "2.4.0"
$ws.Range("C15").Value = 45517
$ws.Range("D15").Value = 45533
$ws.Range("E15").Value = "https://github.com/gaithernOrg/KH1FM-AP/releases"

# Row 16: Kingdom Hearts Re:Chain of Memories
$ws.Range("A16").Value = "Kingdom Hearts Re:Chain of Memories"
$ws.Range("B16").Value = "1.4.2"
$ws.Range("C16").Value = 45516
$ws.Range("D16").Value = 45519
$ws.Range("E16").Value = "https://github.com/gaithernOrg/ArchipelagoKHRECOM/releases"

# Row 17: Metroid Zero Mission
$ws.Range("A17").Value = "Metroid Zero Mission"
$ws.Range("B17").Value = "0.2.0-pre3"
$ws.Range("C17").Value = 45533
$ws.Range("D17").Value = 45533
$ws.Range("E17").Value = "https://github.com/lilDavid/Archipelago-Metroid-Zero-Mission/releases"

# Row 18: Mindustry
$ws.Range("A18").Value = "Mindustry"
$ws.Range("B18").Value = "0.1.1"
$ws.Range("C18").Value = 45528
$ws.Range("D18").Value = 45533
$ws.Range("E18").Value = "https://github.com/JohnMahglass/Archipelago-Mindustry/releases"

# Row 19: Minit
$ws.Range("A19").Value = "Minit"
$ws.Range("B19").Value = "0.6.4"
$ws.Range("C19").Value = 45466
$ws.Range("D19").Value = 45478
$ws.Range("E19").Value = "https://github.com/qwint/APMinit/releases"

# Row 20: Oracle of Seasons
$ws.Range("A20").Value = "Oracle of Seasons"
$ws.Range("B20").Value = "7.2b"
$ws.Range("C20").Value = 45516
$ws.Range("D20").Value = 45519
$ws.Range("E20").Value = "https://github.com/Dinopony/ArchipelagoOoS/releases"

# Row 21: Osu!
$ws.Range("A21").Value = "Osu!"
$ws.Range("B21").Value = "1.0.2"
$ws.Range("C21").Value = 45478
$ws.Range("D21").Value = 45507
$ws.Range("E21").Value = "https://github.com/lilymnky-F/Archipelago-Osu/releases"

# Row 22: Outer Wilds
$ws.Range("A22").Value = "Outer Wilds"
$ws.Range("B22").Value = "0.2.7"
$ws.Range("C22").Value = 45532
$ws.Range("D22").Value = 45533
$ws.Range("E22").Value = "https://github.com/Ixrec/OuterWildsArchipelagoRandomizer/releases"

# Row 23: Pharcryption
$ws.Range("A23").Value = "Pharcryption"
$ws.Range("B23").Value = ""
$ws.Range("C23").Value = 45317
$ws.Range("D23").Value = 45422
$ws.Range("E23").Value = "https://discord.com/channels/731205301247803413/1092478908022136876"

# Row 24: Pokemon Crystal
$ws.Range("A24").Value = "Pokemon Crystal"
$ws.Range("B24").Value = "2.1.2"
$ws.Range("C24").Value = 45488
$ws.Range("D24").Value = 45519
$ws.Range("E24").Value = "https://github.com/AliceMousie/Archipelago/releases"

# Row 25: Pokemon FireRed and LeafGreen
$ws.Range("A25").Value = "Pokemon FireRed and LeafGreen"
$ws.Range("B25").Value = "0.6.0"
$ws.Range("C25").Value = 45533
$ws.Range("D25").Value = 45533
$ws.Range("E25").Value = "https://github.com/vyneras/Archipelago/releases"

# Row 26: Pseudoregalia
$ws.Range("A26").Value = "Pseudoregalia"
$ws.Range("B26").Value = "0.7.2"
$ws.Range("C26").Value = ""
$ws.Range("D26").Value = 45422
$ws.Range("E26").Value = "https://github.com/pseudoregalia-modding/pseudoregalia-archipelago/tags"

# Row 27: Rollercoaster Tycoon 2
$ws.Range("A27").Value = "Rollercoaster Tycoon 2"
$ws.Range("B27").Value = ""
$ws.Range("C27").Value = 45479
$ws.Range("D27").Value = 45533
$ws.Range("E27").Value = "https://github.com/Crazycolbster/rollercoaster-tycoon-randomizer/releases"

# Row 28: Satisfactory
$ws.Range("A28").Value = "Satisfactory"
$ws.Range("B28").Value = "0.1.3.3"
$ws.Range("C28").Value = 45376
$ws.Range("D28").Value = 45533
$ws.Range("E28").Value = "https://github.com/Jarno458/SatisfactoryArchipelagoMod/releases"

# Row 29: Shahrazad
$ws.Range("A29").Value = "Shahrazad"
$ws.Range("B29").Value = "0.1.1"
$ws.Range("C29").Value = 45457
$ws.Range("D29").Value = 45464
$ws.Range("E29").Value = "https://github.com/qwint/ap-shahrazad/releases"

# Row 30: Shapez
$ws.Range("A30").Value = "Shapez"
$ws.Range("B30").Value = "0.3.1"
$ws.Range("C30").Value = 45530
$ws.Range("D30").Value = 45533
$ws.Range("E30").Value = "https://github.com/BlastSlimey/shapezipelago/releases"

# Row 31: Super Metroid Map Rando
$ws.Range("A31").Value = "Super Metroid Map Rando"
$ws.Range("B31").Value = "v111"
$ws.Range("C31").Value = 45442
$ws.Range("D31").Value = 45447
$ws.Range("E31").Value = "https://discord.com/channels/731205301247803413/1156395911874875473"

# Row 32: Super Metroid Subversion
$ws.Range("A32").Value = "Super Metroid Subversion"
$ws.Range("B32").Value = ""
$ws.Range("C32").Value = ""
$ws.Range("D32").Value = 45519
$ws.Range("E32").Value = "https://discord.com/channels/731205301247803413/1146853510378422272"

# Row 33: Wargroove 2
$ws.Range("A33").Value = "Wargroove 2"
$ws.Range("B33").Value = ""
$ws.Range("C33").Value = 45476
$ws.Range("D33").Value = 45478
$ws.Range("E33").Value = "https://discord.com/channels/731205301247803413/1159482310652076082"

# Row 34: Yacht Dice
$ws.Range("A34").Value = "Yacht Dice"
$ws.Range("B34").Value = "2.1.1"
$ws.Range("C34").Value = 45526
$ws.Range("D34").Value = 45533
$ws.Range("E34").Value = "https://github.com/spinerak/ArchipelagoYachtDice/releases"

# --- Re-add hyperlinks to their (possibly new) rows and restore the Hyperlink style ---
$ws.Hyperlinks.Add($ws.Range("E19"), "https://github.com/qwint/APMinit/releases") | Out-Null
$ws.Range("E19").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E16"), "https://github.com/gaithernOrg/ArchipelagoKHRECOM/releases") | Out-Null
$ws.Range("E16").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E15"), "https://github.com/gaithernOrg/KH1FM-AP/releases") | Out-Null
$ws.Range("E15").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E12"), "https://github.com/DrBibop/Archipelago_Inscryption/releases") | Out-Null
$ws.Range("E12").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E9"), "https://github.com/cleartonic/arch_ffvcd/releases") | Out-Null
$ws.Range("E9").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E8"), "https://github.com/Bartz24/Archipelago/releases") | Out-Null
$ws.Range("E8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E6"), "https://github.com/CodeTriangle/CCMultiworldRandomizer/releases") | Out-Null
$ws.Range("E6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E5"), "https://wiki.ctjot.com/doku.php?id=multiworld") | Out-Null
$ws.Range("E5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E26"), "https://github.com/pseudoregalia-modding/pseudoregalia-archipelago/tags") | Out-Null
$ws.Range("E26").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E31"), "https://discord.com/channels/731205301247803413/1156395911874875473") | Out-Null
$ws.Range("E31").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E23"), "https://discord.com/channels/731205301247803413/1092478908022136876") | Out-Null
$ws.Range("E23").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E29"), "https://github.com/qwint/ap-shahrazad/releases") | Out-Null
$ws.Range("E29").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E33"), "https://discord.com/channels/731205301247803413/1159482310652076082") | Out-Null
$ws.Range("E33").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/doshyw/CelesteArchipelago/releases") | Out-Null
$ws.Range("E4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E10"), "https://discord.com/channels/731205301247803413/1022545979146252288") | Out-Null
$ws.Range("E10").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/ScipioWright/Archipelago-SW/releases") | Out-Null
$ws.Range("E2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E22"), "https://github.com/Ixrec/OuterWildsArchipelagoRandomizer/releases") | Out-Null
$ws.Range("E22").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E13"), "https://github.com/ArchipelaGOAL/Archipelago/releases") | Out-Null
$ws.Range("E13").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E11"), "https://github.com/CT075/Archipelago/releases") | Out-Null
$ws.Range("E11").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://discord.com/channels/731205301247803413/1169389087371841708") | Out-Null
$ws.Range("E3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("E7"), "https://github.com/PinkSwitch/Archipelago/releases/") | Out-Null
$ws.Range("E7").Style = "Hyperlink"

# --- Extend the conditional-formatting range to cover the new rows (A2:XFD111) ---
$fc = $ws.Cells.FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("A2:XFD111"))